# plantilla_libros.xlsx edit:
# - Remove the "(Ej: ...)" placeholder example cells from row 2 (columns G-O),
#   leaving G2 with a real sample value ("Madrid") and H2:O2 empty.
# - Add a "Madrid" sample value (sede/lugarPublicacion sample column G) to every
#   data row (3-21), matching the style already used on row 1/2 header cells.
# - Minor view/format touch-ups: header row heights, a few column widths and
#   the saved cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: drop the placeholder "(Ej: ...)" hints in H2:O2, keep the cell
#     style intact (ClearContents preserves formatting, just removes values).
$ws.Range("H2:O2").ClearContents()

# G2 used to hold the placeholder "(Ej: Madrid)" -> becomes a real sample value.
$ws.Range("G2").Value = "Madrid"

# --- Rows 3-21: add a new G column sample value ("Madrid"), copying the
#     formatting (style index) from the already-styled header cell A1 so the
#     new cells pick up the same cellXf the rest of the styled cells use.
$ws.Range("A1").Copy()
$ws.Range("G3:G21").PasteSpecial(-4122)
$ws.Range("G3:G21").Value = "Madrid"

# --- Cosmetic formatting touch-ups ---
# Header rows get an explicit (smaller) row height.
$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(2).RowHeight = 12.75

# A handful of columns get explicit custom widths.
$ws.Columns.Item(1).ColumnWidth = 31
$ws.Columns.Item(3).ColumnWidth = 26.5
$ws.Columns.Item(5).ColumnWidth = 16.833333333333332
$ws.Columns.Item(6).ColumnWidth = 20.666666666666668
$ws.Columns.Item(7).ColumnWidth = 28.166666666666668

# Saved selection moves from H27 to J12.
$ws.Range("J12").Select()

Write-Output "edit applied"
